$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source rows (row 79) used to copy cell formatting (style) for the styled
# columns A (index/bold/border/center style) and E (date-time number format).
$styleSourceA = $ws.Cells.Item(79, 1)
$styleSourceE = $ws.Cells.Item(79, 5)

# --- Row 80 (new match record) ---
$styleSourceA.Copy()
$ws.Cells.Item(80, 1).PasteSpecial(-4122)
$ws.Cells.Item(80, 1).Value = 79

$styleSourceE.Copy()
$ws.Cells.Item(80, 5).PasteSpecial(-4122)
$ws.Cells.Item(80, 5).Value = 45255.5

$ws.Cells.Item(80, 2).Value = "thailand"
$ws.Cells.Item(80, 3).Value = "thai-league-1"
$ws.Cells.Item(80, 4).Value = "2023-2024"
$ws.Cells.Item(80, 6).Value = "Chonburi"
$ws.Cells.Item(80, 8).Value = "Pathum United"
$ws.Cells.Item(80, 11).Value = "18/11/2023 12:12"
$ws.Cells.Item(80, 13).Value = "25/11/2023 11:49"
$ws.Cells.Item(80, 15).Value = "18/11/2023 12:12"
$ws.Cells.Item(80, 17).Value = "25/11/2023 11:53"
$ws.Cells.Item(80, 19).Value = "18/11/2023 12:12"
$ws.Cells.Item(80, 21).Value = "25/11/2023 11:53"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/chonburi-pathum-united/Kh9MlmpI/"
$ws.Cells.Item(80, 7).Value = 1
$ws.Cells.Item(80, 9).Value = 3
$ws.Cells.Item(80, 10).Value = 2.73
$ws.Cells.Item(80, 12).Value = 2.64
$ws.Cells.Item(80, 14).Value = 3.47
$ws.Cells.Item(80, 16).Value = 3.57
$ws.Cells.Item(80, 18).Value = 2.5
$ws.Cells.Item(80, 20).Value = 2.59

# --- Row 81 (new match record) ---
$styleSourceA.Copy()
$ws.Cells.Item(81, 1).PasteSpecial(-4122)
$ws.Cells.Item(81, 1).Value = 80

$styleSourceE.Copy()
$ws.Cells.Item(81, 5).PasteSpecial(-4122)
$ws.Cells.Item(81, 5).Value = 45255.54166666666

$ws.Cells.Item(81, 2).Value = "thailand"
$ws.Cells.Item(81, 3).Value = "thai-league-1"
$ws.Cells.Item(81, 4).Value = "2023-2024"
$ws.Cells.Item(81, 6).Value = "Sukhothai"
$ws.Cells.Item(81, 8).Value = "Prachuap"
$ws.Cells.Item(81, 11).Value = "18/11/2023 13:13"
$ws.Cells.Item(81, 13).Value = "25/11/2023 11:20"
$ws.Cells.Item(81, 15).Value = "18/11/2023 13:13"
$ws.Cells.Item(81, 17).Value = "25/11/2023 11:20"
$ws.Cells.Item(81, 19).Value = "18/11/2023 13:13"
$ws.Cells.Item(81, 21).Value = "25/11/2023 12:51"
$ws.Cells.Item(81, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/sukhothai-prachuap/t2AIkTUB/"
$ws.Cells.Item(81, 7).Value = 2
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 2.32
$ws.Cells.Item(81, 12).Value = 2.37
$ws.Cells.Item(81, 14).Value = 3.53
$ws.Cells.Item(81, 16).Value = 3.65
$ws.Cells.Item(81, 18).Value = 2.95
$ws.Cells.Item(81, 20).Value = 2.87

# --- Row 82 (new match record) ---
$styleSourceA.Copy()
$ws.Cells.Item(82, 1).PasteSpecial(-4122)
$ws.Cells.Item(82, 1).Value = 81

$styleSourceE.Copy()
$ws.Cells.Item(82, 5).PasteSpecial(-4122)
$ws.Cells.Item(82, 5).Value = 45255.58333333334

$ws.Cells.Item(82, 2).Value = "thailand"
$ws.Cells.Item(82, 3).Value = "thai-league-1"
$ws.Cells.Item(82, 4).Value = "2023-2024"
$ws.Cells.Item(82, 6).Value = "Lamphun Warrior"
$ws.Cells.Item(82, 8).Value = "Muang Thong Utd"
$ws.Cells.Item(82, 11).Value = "18/11/2023 14:12"
$ws.Cells.Item(82, 13).Value = "25/11/2023 13:51"
$ws.Cells.Item(82, 15).Value = "18/11/2023 14:12"
$ws.Cells.Item(82, 17).Value = "25/11/2023 13:51"
$ws.Cells.Item(82, 19).Value = "18/11/2023 14:12"
$ws.Cells.Item(82, 21).Value = "25/11/2023 13:51"
$ws.Cells.Item(82, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/lamphun-warrior-muang-thong-utd/bwLDj9F5/"
$ws.Cells.Item(82, 7).Value = 3
$ws.Cells.Item(82, 9).Value = 1
$ws.Cells.Item(82, 10).Value = 2.69
$ws.Cells.Item(82, 12).Value = 2.3
$ws.Cells.Item(82, 14).Value = 3.44
$ws.Cells.Item(82, 16).Value = 3.57
$ws.Cells.Item(82, 18).Value = 2.55
$ws.Cells.Item(82, 20).Value = 3.05

# --- Row 83 (new match record) ---
$styleSourceA.Copy()
$ws.Cells.Item(83, 1).PasteSpecial(-4122)
$ws.Cells.Item(83, 1).Value = 82

$styleSourceE.Copy()
$ws.Cells.Item(83, 5).PasteSpecial(-4122)
$ws.Cells.Item(83, 5).Value = 45256.47916666666

$ws.Cells.Item(83, 2).Value = "thailand"
$ws.Cells.Item(83, 3).Value = "thai-league-1"
$ws.Cells.Item(83, 4).Value = "2023-2024"
$ws.Cells.Item(83, 6).Value = "Uthai Thani"
$ws.Cells.Item(83, 8).Value = "Police Tero"
$ws.Cells.Item(83, 11).Value = "19/11/2023 11:43"
$ws.Cells.Item(83, 13).Value = "26/11/2023 11:02"
$ws.Cells.Item(83, 15).Value = "19/11/2023 11:43"
$ws.Cells.Item(83, 17).Value = "26/11/2023 11:02"
$ws.Cells.Item(83, 19).Value = "19/11/2023 11:43"
$ws.Cells.Item(83, 21).Value = "26/11/2023 11:02"
$ws.Cells.Item(83, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/uthai-thani-police-tero/hKMTnREU/"
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 9).Value = 1
$ws.Cells.Item(83, 10).Value = 1.71
$ws.Cells.Item(83, 12).Value = 1.97
$ws.Cells.Item(83, 14).Value = 4.11
$ws.Cells.Item(83, 16).Value = 3.92
$ws.Cells.Item(83, 18).Value = 4.01
$ws.Cells.Item(83, 20).Value = 3.52

# --- Row 84 (new match record) ---
$styleSourceA.Copy()
$ws.Cells.Item(84, 1).PasteSpecial(-4122)
$ws.Cells.Item(84, 1).Value = 83

$styleSourceE.Copy()
$ws.Cells.Item(84, 5).PasteSpecial(-4122)
$ws.Cells.Item(84, 5).Value = 45256.54166666666

$ws.Cells.Item(84, 2).Value = "thailand"
$ws.Cells.Item(84, 3).Value = "thai-league-1"
$ws.Cells.Item(84, 4).Value = "2023-2024"
$ws.Cells.Item(84, 6).Value = "Chiangrai Utd"
$ws.Cells.Item(84, 8).Value = "Port MTI FC"
$ws.Cells.Item(84, 11).Value = "19/11/2023 13:12"
$ws.Cells.Item(84, 13).Value = "26/11/2023 12:51"
$ws.Cells.Item(84, 15).Value = "19/11/2023 13:12"
$ws.Cells.Item(84, 17).Value = "26/11/2023 12:56"
$ws.Cells.Item(84, 19).Value = "19/11/2023 13:12"
$ws.Cells.Item(84, 21).Value = "26/11/2023 12:56"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/chiangrai-utd-port-mti-fc/pzH9ik0a/"
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(84, 10).Value = 4.61
$ws.Cells.Item(84, 12).Value = 4.16
$ws.Cells.Item(84, 14).Value = 4.03
$ws.Cells.Item(84, 16).Value = 3.76
$ws.Cells.Item(84, 18).Value = 1.63
$ws.Cells.Item(84, 20).Value = 1.85

# Clear clipboard marching ants / leftover clipboard state.
$excel.CutCopyMode = 0